$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I3 value changes 44 -> 46
$ws.Range("I3").Value = 46

# New TIME series used values for GHB (row25), RIV (row26), WEL (row31), DRN (row32), RCH (row33)
$ws.Range("I25").Value = 0.75
$ws.Range("I26").Value = 0.75
$ws.Range("I31").Value = 0.75
$ws.Range("I32").Value = 2
$ws.Range("I33").Value = 1

# Row 31 height set to 15 (custom height)
$ws.Rows.Item(31).RowHeight = 15

# Selection changes to I33
$ws.Range("I33").Select()
